$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.197736333333333
$ws.Range("N2").Value = 9.593208999999998
$ws.Range("Q2").Value = 0.01296682083166667
$ws.Range("R2").Value = 0.116701387485
